$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 0
$ws.Range("O6").Value = 0.001269340515136719
$ws.Range("O11").Value = 0.03039073944091797
$ws.Range("O12").Value = 0.0157017707824707
$ws.Range("O13").Value = 0.01247143745422363
$ws.Range("O14").Value = 0.03470206260681152
$ws.Range("O15").Value = 0.002999782562255859
$ws.Range("O17").Value = 0.001062393188476562
$ws.Range("O18").Value = 0.03235220909118652
$ws.Range("O19").Value = 0.02439498901367188
$ws.Range("O21").Value = 0.001807928085327148
$ws.Range("O24").Value = 0.0006709098815917969
$ws.Range("O25").Value = 0.0319511890411377
$ws.Range("O26").Value = 0.03786134719848633
$ws.Range("O28").Value = 0.01002025604248047
$ws.Range("O30").Value = 0.002600193023681641
$ws.Range("O33").Value = 0
$ws.Range("O34").Value = 6.898502588272095
$ws.Range("O35").Value = 0.108687162399292
$ws.Range("O36").Value = 0.002005100250244141
$ws.Range("O37").Value = 0.8772985935211182
$ws.Range("O38").Value = 0.1683557033538818
$ws.Range("O39").Value = 0.1531956195831299
$ws.Range("O40").Value = 0.003763198852539062
$ws.Range("O41").Value = 0.005013704299926758
$ws.Range("O42").Value = 0.003263473510742188
$ws.Range("O43").Value = 0.1390199661254883
$ws.Range("O44").Value = 0.1347277164459229
$ws.Range("O45").Value = 0.1537575721740723
$ws.Range("O46").Value = 0.1164121627807617
$ws.Range("O47").Value = 0.1105556488037109
$ws.Range("O48").Value = 0.08990836143493652
$ws.Range("O49").Value = 0.02728581428527832
$ws.Range("O50").Value = 0.002504825592041016
$ws.Range("O51").Value = 0.001003026962280273
$ws.Range("O52").Value = 0.0659325122833252
$ws.Range("O53").Value = 0
$ws.Range("O55").Value = 0.004694938659667969
$ws.Range("O57").Value = 0.01651144027709961
$ws.Range("O58").Value = 0.01616024971008301
$ws.Range("O59").Value = 0.06256008148193359
$ws.Range("O60").Value = 0.04228949546813965
$ws.Range("O61").Value = 0
$ws.Range("O62").Value = 0.03123855590820312
$ws.Range("O63").Value = 0.004289627075195312
$ws.Range("O64").Value = 0.05036377906799316
$ws.Range("O65").Value = 0.262686014175415
$ws.Range("O66").Value = 0.0741569995880127
$ws.Range("O67").Value = 0.01680517196655273
$ws.Range("O68").Value = 0
$ws.Range("O69").Value = 0.0004804134368896484
$ws.Range("O71").Value = 11.95888018608093
$ws.Range("O72").Value = 25.44644594192505
$ws.Range("O73").Value = 0.2148764133453369
$ws.Range("O74").Value = 0.1177046298980713
$ws.Range("O75").Value = 0.0665290355682373
$ws.Range("O76").Value = 0.08156180381774902
$ws.Range("O77").Value = 0.004179239273071289
$ws.Range("O79").Value = 0.04972338676452637
$ws.Range("O80").Value = 0.007565021514892578
$ws.Range("O81").Value = 0.3901727199554443
$ws.Range("O82").Value = 0.1012170314788818
$ws.Range("O83").Value = 0.09713125228881836
$ws.Range("O84").Value = 1.490601539611816
$ws.Range("O85").Value = 0.06674647331237793
$ws.Range("O87").Value = 0.06664490699768066
$ws.Range("O88").Value = 0.108468770980835
$ws.Range("O89").Value = 0.09443116188049316
$ws.Range("O90").Value = 8.243488788604736
$ws.Range("O91").Value = 0.05160951614379883
$ws.Range("O92").Value = 0.004177331924438477
$ws.Range("O93").Value = 0.02486920356750488
$ws.Range("O94").Value = 0
$ws.Range("O95").Value = 0.01563715934753418
$ws.Range("O96").Value = 0.3589973449707031
$ws.Range("O97").Value = 0.1582858562469482
$ws.Range("O98").Value = 0.09331226348876953
$ws.Range("O99").Value = 0.09979391098022461
$ws.Range("O100").Value = 0.07330727577209473
$ws.Range("O101").Value = 0
$ws.Range("O103").Value = 0.009509563446044922
$ws.Range("O104").Value = 17.05698871612549
$ws.Range("O105").Value = 0.06062889099121094
$ws.Range("O106").Value = 0.06898951530456543
$ws.Range("O107").Value = 0.1340494155883789
$ws.Range("O108").Value = 0.0006802082061767578
$ws.Range("O109").Value = 0
$ws.Range("O110").Value = 0.02241897583007812
$ws.Range("O111").Value = 0.04070329666137695
$ws.Range("O112").Value = 0.07308864593505859
$ws.Range("O113").Value = 0.02842450141906738
$ws.Range("O114").Value = 0.02129507064819336
$ws.Range("O115").Value = 0.0165705680847168
$ws.Range("O116").Value = 0.02194476127624512
$ws.Range("O117").Value = 0.03135347366333008
$ws.Range("O121").Value = 0.001670598983764648
$ws.Range("O122").Value = 0.07932186126708984
$ws.Range("O125").Value = 0
$ws.Range("O126").Value = 0
$ws.Range("O127").Value = 0.001067876815795898
$ws.Range("O129").Value = 0
$ws.Range("O130").Value = 0.001007556915283203
$ws.Range("O132").Value = 0.003516912460327148
$ws.Range("O133").Value = 0
$ws.Range("O134").Value = 0.001000165939331055
$ws.Range("O135").Value = 0
$ws.Range("O136").Value = 0.01550555229187012
$ws.Range("O137").Value = 0.006011724472045898
$ws.Range("O138").Value = 0
$ws.Range("O139").Value = 0.008244991302490234
$ws.Range("O142").Value = 0.01000785827636719
$ws.Range("O143").Value = 0.00153040885925293
$ws.Range("O144").Value = 0.01700067520141602
$ws.Range("O145").Value = 0
$ws.Range("O148").Value = 0.03065800666809082
$ws.Range("O149").Value = 0.1173393726348877
$ws.Range("O150").Value = 1.26715612411499
$ws.Range("O151").Value = 0.0534825325012207
$ws.Range("O152").Value = 0.0003843307495117188
$ws.Range("O153").Value = 0.0009455680847167969
$ws.Range("O154").Value = 0.002998590469360352
$ws.Range("O155").Value = 0
$ws.Range("O156").Value = 0.0115058422088623
$ws.Range("O157").Value = 0
$ws.Range("O159").Value = 0.03674674034118652
$ws.Range("O160").Value = 0
$ws.Range("O161").Value = 0.01563405990600586
$ws.Range("O162").Value = 0.006996631622314453
$ws.Range("O164").Value = 0.04327082633972168
$ws.Range("O165").Value = 0
$ws.Range("O166").Value = 0.0218353271484375
$ws.Range("O167").Value = 0.06928706169128418
$ws.Range("O169").Value = 0
$ws.Range("O170").Value = 0.03306221961975098
$ws.Range("O171").Value = 0.01563572883605957
